$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking.com scrape refresh (GitHub Actions cron).
#
# Price (col D) and Volume/1h (col E) are stored as literal TEXT in the source sheet
# (e.g. "27.453.80", "1.005", "  -0.72%  ") even though several Price strings parse as
# plain numbers ("1.005", "332.06", ...). Assigning such a string straight to .Value
# makes Excel auto-convert it to a real number, which would change both the stored type
# and silently normalise the text (e.g. "156.60" -> 156.6). To keep those cells text,
# round-trip them through a text formula + paste-as-values, which stores a literal string
# without touching the cell number format/style.
function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $escaped = $val -replace '"', '""'
    $c.Formula = "=""" + $escaped + """"
    $c.Copy()
    $c.PasteSpecial(-4163)
}

$ws.Range("D2").Value = "27.453.80"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "1.824.22"
$ws.Range("E3").Value = "  -2.08%  "
Set-TextValue "D4" "1.005"
$ws.Range("E4").Value = "  -0.78%  "
Set-TextValue "D5" "332.06"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  -0.71%  "
Set-TextValue "D7" "0.4541"
$ws.Range("E7").Value = "  -2.20%  "
Set-TextValue "D8" "0.3801"
$ws.Range("E8").Value = "  -2.18%  "
Set-TextValue "D9" "46.23"
$ws.Range("E9").Value = "  +0.32%  "
Set-TextValue "D10" "0.07869"
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("E11").Value = "  -2.58%  "
Set-TextValue "D12" "20.97"
$ws.Range("E12").Value = "  -2.49%  "
$ws.Range("D13").Value = "1.825.87"
$ws.Range("E13").Value = "  -1.73%  "
Set-TextValue "D14" "5.858"
$ws.Range("E14").Value = "  -1.91%  "
Set-TextValue "D15" "7.026"
$ws.Range("E15").Value = "  -2.06%  "
Set-TextValue "D16" "1.008"
$ws.Range("E16").Value = "  -0.50%  "
Set-TextValue "D17" "88.51"
$ws.Range("E17").Value = "  +0.79%  "
Set-TextValue "D18" "0.06629"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("E19").Value = "  -1.46%  "
Set-TextValue "D20" "17.16"
$ws.Range("E20").Value = "  +1.67%  "
Set-TextValue "D21" "1.004"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "27.428.74"
$ws.Range("E22").Value = "  -0.76%  "
Set-TextValue "D23" "5.317"
$ws.Range("E23").Value = "  -2.39%  "
Set-TextValue "D24" "10.77"
$ws.Range("E24").Value = "  -0.73%  "
Set-TextValue "D25" "2.302"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").Value = "2.050.20"
$ws.Range("E26").Value = "  -1.49%  "
Set-TextValue "D27" "156.60"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("E28").Value = "  -1.64%  "
Set-TextValue "D29" "2.056"
$ws.Range("E29").Value = "  -2.31%  "
Set-TextValue "D30" "5.237"
$ws.Range("E30").Value = "  -2.13%  "
Set-TextValue "D31" "117.88"
$ws.Range("E31").Value = "  -2.81%  "
Set-TextValue "D32" "0.9439"
$ws.Range("E32").Value = "  -2.57%  "
Set-TextValue "D33" "0.09294"
$ws.Range("E33").Value = "  -1.56%  "
Set-TextValue "D34" "3.576"
$ws.Range("E34").Value = "  -1.92%  "
Set-TextValue "D35" "5.226"
$ws.Range("E35").Value = "  -1.10%  "
Set-TextValue "D36" "1.317"
$ws.Range("E36").Value = "  +0.54%  "
Set-TextValue "D37" "0.05911"
$ws.Range("E37").Value = "  -1.67%  "
Set-TextValue "D38" "0.02184"
$ws.Range("E38").Value = "  -1.36%  "
Set-TextValue "D39" "1.155"
$ws.Range("E39").Value = "  -3.58%  "
Set-TextValue "D40" "8.013"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("E41").Value = "  -2.62%  "
Set-TextValue "D42" "0.1828"
$ws.Range("E42").Value = "  -2.56%  "
Set-TextValue "D43" "9.996"
$ws.Range("E43").Value = "  -2.19%  "
Set-TextValue "D44" "1.276"
$ws.Range("E44").Value = "  +1.87%  "
Set-TextValue "D47" "1.861"
$ws.Range("E47").Value = "  -2.68%  "
Set-TextValue "D48" "0.06605"
$ws.Range("E48").Value = "  -2.26%  "
Set-TextValue "D49" "110.22"
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("E51").Value = "  -0.83%  "

# Rows 45 and 46 swapped order: EnergySwap now ranks above Decentraland
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "11.95"
$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.5440"
$ws.Range("E46").Value = "  -2.81%  "

$excel.CutCopyMode = $false
